$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update language code column (A) from French to English
$ws.Range("A2").Value2 = "eng"
$ws.Range("A3").Value2 = "eng"
$ws.Range("A4").Value2 = "eng"
$ws.Range("A5").Value2 = "eng"
$ws.Range("A6").Value2 = "eng"
$ws.Range("A7").Value2 = "eng"

# Update name/description text columns (C/D) with English translations
$ws.Range("C2").Value2 = "Pre-Registration"
$ws.Range("D2").Value2 = "Web portal for pre-registrations"

$ws.Range("C3").Value2 = "Registration Client"
$ws.Range("D3").Value2 = "Desktop application for Registrations"

$ws.Range("C4").Value2 = "Registration Processor"
$ws.Range("D4").Value2 = "Application for post-registration process"

$ws.Range("C5").Value2 = "ID Authentication"
$ws.Range("D5").Value2 = "Application for third party service provider authentication"

$ws.Range("C6").Value2 = "ID Control"
$ws.Range("D6").Value2 = "Web portal for configuring applications"

$ws.Range("C7").Value2 = "Resident Portal"
$ws.Range("D7").Value2 = "Web portal for Post ID generation services"

# Row heights grew to accommodate the longer English text wrapping
$ws.Rows.Item(2).RowHeight = 66
$ws.Rows.Item(3).RowHeight = 66
$ws.Rows.Item(4).RowHeight = 79
$ws.Rows.Item(5).RowHeight = 105
$ws.Rows.Item(6).RowHeight = 79
$ws.Rows.Item(7).RowHeight = 66

# Selection moved
$ws.Range("H3").Select()
